$wb = $excel.ActiveWorkbook

# --- Controls sheet: update n_years value, begin adding new rows ---
$ws1 = $wb.Worksheets.Item("Controls")
$ws1.Range("B2").Value = 200
$ws1.Range("A5").Value = "n_sex"
$ws1.Range("B5").Value = 1

# --- Maturity_At_Age: insert a Sex column, add sex=2 row ---
$ws3 = $wb.Worksheets.Item("Maturity_At_Age")
$ws3.Columns("C:C").Insert()
$ws3.Range("C1").Value = "Sex"
$ws3.Range("C2").Value = 1
$ws3.Range("A3").Value = 1
$ws3.Range("B3").Value = "Time_Inv"
$ws3.Range("C3").Value = 2
$ws3.Cells.Item(3, 4).Value = 0.022375363
$ws3.Cells.Item(3, 5).Value = 0.046018622
$ws3.Cells.Item(3, 6).Value = 0.092286518
$ws3.Cells.Item(3, 7).Value = 0.176467916
$ws3.Cells.Item(3, 8).Value = 0.311118342
$ws3.Cells.Item(3, 9).Value = 0.487670289
$ws3.Cells.Item(3, 10).Value = 0.667353261
$ws3.Cells.Item(3, 11).Value = 0.808734125
$ws3.Cells.Item(3, 12).Value = 0.899109861
$ws3.Cells.Item(3, 13).Value = 0.949450979
$ws3.Cells.Item(3, 14).Value = 0.975361733
$ws3.Cells.Item(3, 15).Value = 0.988156652
$ws3.Cells.Item(3, 16).Value = 0.994345552
$ws3.Cells.Item(3, 17).Value = 0.997309166
$ws3.Cells.Item(3, 18).Value = 0.998721485
$ws3.Cells.Item(3, 19).Value = 0.999392981
$ws3.Cells.Item(3, 20).Value = 0.999711899
$ws3.Cells.Item(3, 21).Value = 0.999863285
$ws3.Cells.Item(3, 22).Value = 0.999935129
$ws3.Cells.Item(3, 23).Value = 0.99996922
$ws3.Cells.Item(3, 24).Value = 0.999985396
$ws3.Cells.Item(3, 25).Value = 0.999993071
$ws3.Cells.Item(3, 26).Value = 0.999996712
$ws3.Cells.Item(3, 27).Value = 0.99999844
$ws3.Cells.Item(3, 28).Value = 0.99999926
$ws3.Cells.Item(3, 29).Value = 0.999999649
$ws3.Cells.Item(3, 30).Value = 0.999999833
$ws3.Cells.Item(3, 31).Value = 0.999999921
$ws3.Cells.Item(3, 32).Value = 0.999999962
$ws3.Cells.Item(3, 33).Value = 0.999999982

# --- Weight_At_Age: insert a Sex column, add sex=2 row with new weights ---
$ws4 = $wb.Worksheets.Item("Weight_At_Age")
$ws4.Columns("C:C").Insert()
$ws4.Range("C1").Value = "Sex"
$ws4.Range("C2").Value = 1
$ws4.Range("A3").Value = 1
$ws4.Range("B3").Value = "Time_Inv"
$ws4.Range("C3").Value = 2
$ws4.Cells.Item(3, 4).Value = 1.1085
$ws4.Cells.Item(3, 5).Value = 1.4285
$ws4.Cells.Item(3, 6).Value = 1.7228
$ws4.Cells.Item(3, 7).Value = 1.9837
$ws4.Cells.Item(3, 8).Value = 2.2089
$ws4.Cells.Item(3, 9).Value = 2.3995
$ws4.Cells.Item(3, 10).Value = 2.5586
$ws4.Cells.Item(3, 11).Value = 2.6899
$ws4.Cells.Item(3, 12).Value = 2.7974
$ws4.Cells.Item(3, 13).Value = 2.8848
$ws4.Cells.Item(3, 14).Value = 2.9555
$ws4.Cells.Item(3, 15).Value = 3.0125
$ws4.Cells.Item(3, 16).Value = 3.0584
$ws4.Cells.Item(3, 17).Value = 3.0951
$ws4.Cells.Item(3, 18).Value = 3.1245
$ws4.Cells.Item(3, 19).Value = 3.148
$ws4.Cells.Item(3, 20).Value = 3.1668
$ws4.Cells.Item(3, 21).Value = 3.1817
$ws4.Cells.Item(3, 22).Value = 3.1936
$ws4.Cells.Item(3, 23).Value = 3.2031
$ws4.Cells.Item(3, 24).Value = 3.2107
$ws4.Cells.Item(3, 25).Value = 3.2167
$ws4.Cells.Item(3, 26).Value = 3.2215
$ws4.Cells.Item(3, 27).Value = 3.2253
$ws4.Cells.Item(3, 28).Value = 3.2283
$ws4.Cells.Item(3, 29).Value = 3.2307
$ws4.Cells.Item(3, 30).Value = 3.2326
$ws4.Cells.Item(3, 31).Value = 3.2341
$ws4.Cells.Item(3, 32).Value = 3.2353
$ws4.Cells.Item(3, 33).Value = 3.2381

# --- Back to Controls: finish adding new rows 5 and 6 ---
$ws1.Range("C5").Value = "Number of sexes (1. = Females, 2 = Males)"
$ws1.Range("A6").Value = "n_fleets"
$ws1.Range("B6").Value = 1
$ws1.Range("C6").Value = "Number of fishery fleets"

# --- Recruitment_Mortality: update sigma_rec value ---
$ws5 = $wb.Worksheets.Item("Recruitment_Mortality")
$ws5.Range("B5").Value = 0.5

# --- View / selection state ---
$ws5.Range("B5").Select()
$ws4.Range("C2").Select()
$ws3.Range("D3:AG3").Select()
$ws1.Activate()
$ws1.Range("C12").Select()
